$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.994.90"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.16"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.61"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4643"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3870"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07854"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9579"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.81"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.816.34"
$ws.Range("E12").Value = "  -7.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.666"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.893"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06779"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.21"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009904"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.58"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.009.75"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.311"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.091"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.047.12"
$ws.Range("E25").Value = "  -6.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.65"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.733"
$ws.Range("E28").Value = "  -7.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.969"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.22"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09253"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9338"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.284"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.317"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.287"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05856"
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02141"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.143"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.777"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5574"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.847"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1760"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5250"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06997"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.116"
$ws.Range("E46").Value = "  -11.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.110"
$ws.Range("E47").Value = "  -11.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.826"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.84"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.324"
$ws.Range("E51").Value = "  +0.58%  "
